# HLA_Attribute_Classification.xlsx - "all attributed added to excel"
#
# This script:
#  1. Removes the yellow highlight fill from the first few attribute cells
#     in column A (rows 2-4, 6-11, 13-20) that no longer need emphasis.
#  2. Appends two new attribute rows at the bottom of the table:
#       A62: tce_imm_match
#       A63: cyto_score_detail
#  3. Updates the active view/selection to reflect where the user ended up
#     working (scrolled down near the bottom of the sheet, selection on F61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the yellow fill/style from the previously-highlighted cells ---
$highlightedRows = @(2,3,4,6,7,8,9,10,11,13,14,15,16,17,18,19,20)
foreach ($r in $highlightedRows) {
    $ws.Cells.Item($r, 1).ClearFormats()
}

# --- 2. Add the two newly-identified attributes at the end of the table ---
$ws.Cells.Item(62, 1).Value = "tce_imm_match"
$ws.Cells.Item(63, 1).Value = "cyto_score_detail"

# --- 3. Reflect the final scroll position / selection used while editing ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("F61").Select()
